# Update "想去人数" (column F) values on the "展览" sheet and the
# corresponding rows on the "全部类型" sheet (which aggregates all
# sheets, offset by one row versus "展览").

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll = $wb.Worksheets.Item("全部类型")

# Row => New value, for the "展览" sheet
$exhibitUpdates = @{
    2  = 1196
    4  = 638
    7  = 59
    9  = 361
    11 = 114
    13 = 279
    18 = 288
    19 = 720
    20 = 101
    21 = 682
    22 = 220
    23 = 53
    25 = 389
    26 = 210
    29 = 19
    30 = 27
    31 = 438
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# Row => New value, for the "全部类型" sheet (same rows shifted by +1)
$allUpdates = @{
    3  = 1196
    5  = 638
    9  = 59
    11 = 361
    13 = 114
    15 = 279
    25 = 288
    26 = 720
    27 = 101
    28 = 682
    29 = 220
    30 = 53
    32 = 389
    35 = 210
    40 = 19
    41 = 27
    43 = 438
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}
